$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 125
$ws.Range("F4").Value = 643
$ws.Range("F5").Value = 371
$ws.Range("F6").Value = 550
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 11760
$ws.Range("F13").Value = 2112
$ws.Range("F15").Value = 238
$ws.Range("F18").Value = 1198
$ws.Range("F20").Value = 251
$ws.Range("F21").Value = 743
$ws.Range("F22").Value = 659
$ws.Range("F23").Value = 274
$ws.Range("F24").Value = 2403
$ws.Range("F25").Value = 731
$ws.Range("F26").Value = 3659
$ws.Range("F27").Value = 3659
$ws.Range("F28").Value = 1069
$ws.Range("F29").Value = 818
$ws.Range("F33").Value = 994
$ws.Range("F35").Value = 72
$ws.Range("F36").Value = 256
$ws.Range("F37").Value = 21
$ws.Range("F39").Value = 15
$ws.Range("F40").Value = 3351
$ws.Range("F41").Value = 4444
$ws.Range("F42").Value = 5484
$ws.Range("F43").Value = 106
$ws.Range("F45").Value = 158
$ws.Range("F46").Value = 273
$ws.Range("F48").Value = 32
$ws.Range("F49").Value = 4093
$ws.Range("F50").Value = 104

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4157
$ws.Range("F12").Value = 731

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 65

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 65
$ws.Range("F6").Value = 125
$ws.Range("F7").Value = 643
$ws.Range("F8").Value = 371
$ws.Range("F9").Value = 550
$ws.Range("F11").Value = 11760
$ws.Range("F15").Value = 2112
$ws.Range("F18").Value = 1198
$ws.Range("F20").Value = 251
$ws.Range("F21").Value = 4157
$ws.Range("F23").Value = 274
$ws.Range("F24").Value = 731
$ws.Range("F25").Value = 3659
$ws.Range("F26").Value = 1069
$ws.Range("F29").Value = 818
$ws.Range("F31").Value = 994
$ws.Range("F33").Value = 72
$ws.Range("F34").Value = 256
$ws.Range("F35").Value = 21
$ws.Range("F36").Value = 15
$ws.Range("F37").Value = 4444
$ws.Range("F38").Value = 106
$ws.Range("F40").Value = 158
$ws.Range("F41").Value = 273
$ws.Range("F45").Value = 32
$ws.Range("F50").Value = 104
